# Scheduled-runner market price refresh for Ultros_Profits workbook.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ,
# LevePriceNQ / LevePriceHQ and the derived LeveProfitNQ / LeveProfitHQ columns
# (H, I, J, K, L, M, N) for the specific leve rows whose market data changed.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
# Row 40
$ws_ALC.Range("H40").Value = 2912.8333
$ws_ALC.Range("I40").Value = 1977.5
$ws_ALC.Range("K40").Value = 1977.5
$ws_ALC.Range("M40").Value = -1802.5
# Row 62
$ws_ALC.Range("H62").Value = 2832.5
$ws_ALC.Range("I62").Value = 2399
$ws_ALC.Range("K62").Value = 2399
$ws_ALC.Range("M62").Value = -1775
# Row 65
$ws_ALC.Range("H65").Value = 2832.5
$ws_ALC.Range("I65").Value = 2399
$ws_ALC.Range("K65").Value = 11995
$ws_ALC.Range("M65").Value = -8875
# Row 76
$ws_ALC.Range("H76").Value = 5220
$ws_ALC.Range("I76").Value = 5375
$ws_ALC.Range("K76").Value = 5375
$ws_ALC.Range("M76").Value = -5060
# Row 79
$ws_ALC.Range("H79").Value = 5220
$ws_ALC.Range("I79").Value = 5375
$ws_ALC.Range("K79").Value = 5375
$ws_ALC.Range("M79").Value = -4283
# Row 98
$ws_ALC.Range("H98").Value = 2204.2
$ws_ALC.Range("I98").Value = 1966.3846
$ws_ALC.Range("K98").Value = 1966.3846
$ws_ALC.Range("M98").Value = -468.3846000000001
# Row 122
$ws_ALC.Range("H122").Value = 2204.2
$ws_ALC.Range("I122").Value = 1966.3846
$ws_ALC.Range("K122").Value = 5899.1538
$ws_ALC.Range("M122").Value = -3449.1538
# Row 132
$ws_ALC.Range("H132").Value = 14123.975
$ws_ALC.Range("I132").Value = 997.89655
$ws_ALC.Range("K132").Value = 2993.68965
$ws_ALC.Range("M132").Value = -463.6896500000003
# Row 135
$ws_ALC.Range("H135").Value = 3048.0557
$ws_ALC.Range("I135").Value = 3295.875
$ws_ALC.Range("J135").Value = 1065.5
$ws_ALC.Range("K135").Value = 29662.875
$ws_ALC.Range("L135").Value = 9589.5
$ws_ALC.Range("M135").Value = -27127.875
$ws_ALC.Range("N135").Value = -14659.5
# --- ARM ---
# Row 61
$ws_ARM.Range("H61").Value = 3286.923
$ws_ARM.Range("I61").Value = 1498.1666
$ws_ARM.Range("J61").Value = 7311.625
$ws_ARM.Range("K61").Value = 1498.1666
$ws_ARM.Range("L61").Value = 7311.625
$ws_ARM.Range("M61").Value = -1286.1666
$ws_ARM.Range("N61").Value = -7735.625
# Row 122
$ws_ARM.Range("H122").Value = 3463.75
$ws_ARM.Range("I122").Value = 2142.3333
$ws_ARM.Range("J122").Value = 5445.875
$ws_ARM.Range("K122").Value = 6426.999899999999
$ws_ARM.Range("L122").Value = 16337.625
$ws_ARM.Range("M122").Value = -3976.999899999999
$ws_ARM.Range("N122").Value = -21237.625
# Row 132
$ws_ARM.Range("H132").Value = 550.5714
$ws_ARM.Range("I132").Value = 550.5714
$ws_ARM.Range("K132").Value = 1651.7142
$ws_ARM.Range("M132").Value = 878.2857999999999
# Row 136
$ws_ARM.Range("H136").Value = 3286.923
$ws_ARM.Range("I136").Value = 1498.1666
$ws_ARM.Range("J136").Value = 7311.625
$ws_ARM.Range("K136").Value = 4494.4998
$ws_ARM.Range("L136").Value = 21934.875
$ws_ARM.Range("M136").Value = -1944.4998
$ws_ARM.Range("N136").Value = -27034.875
# --- BSM ---
# Row 20
$ws_BSM.Range("H20").Value = 3447
$ws_BSM.Range("I20").Value = 2195.111
$ws_BSM.Range("J20").Value = 6666.143
$ws_BSM.Range("K20").Value = 2195.111
$ws_BSM.Range("L20").Value = 6666.143
$ws_BSM.Range("M20").Value = -1948.111
$ws_BSM.Range("N20").Value = -7160.143
# Row 94
$ws_BSM.Range("H94").Value = 5267.933
$ws_BSM.Range("J94").Value = 5971.4287
$ws_BSM.Range("L94").Value = 5971.4287
$ws_BSM.Range("N94").Value = -6873.4287
# Row 134
$ws_BSM.Range("H134").Value = 3311.6
$ws_BSM.Range("I134").Value = 1915.8334
$ws_BSM.Range("K134").Value = 5747.5002
$ws_BSM.Range("M134").Value = -3212.5002
# --- CRP ---
# Row 99
$ws_CRP.Range("H99").Value = 11816113
$ws_CRP.Range("I99").Value = 2714793.5
$ws_CRP.Range("J99").Value = 20007300
$ws_CRP.Range("K99").Value = 2714793.5
$ws_CRP.Range("L99").Value = 20007300
$ws_CRP.Range("M99").Value = -2713295.5
$ws_CRP.Range("N99").Value = -20010296
# Row 126
$ws_CRP.Range("H126").Value = 11816113
$ws_CRP.Range("I126").Value = 2714793.5
$ws_CRP.Range("J126").Value = 20007300
$ws_CRP.Range("K126").Value = 8144380.5
$ws_CRP.Range("L126").Value = 60021900
$ws_CRP.Range("M126").Value = -8141910.5
$ws_CRP.Range("N126").Value = -60026840
# Row 134
$ws_CRP.Range("H134").Value = 2897.5483
$ws_CRP.Range("I134").Value = 1678.9333
$ws_CRP.Range("K134").Value = 5036.7999
$ws_CRP.Range("M134").Value = -2501.7999
# --- CUL ---
# Row 122
$ws_CUL.Range("H122").Value = 1760.8889
$ws_CUL.Range("I122").Value = 1063.5
$ws_CUL.Range("J122").Value = 2318.8
$ws_CUL.Range("K122").Value = 9571.5
$ws_CUL.Range("L122").Value = 20869.2
$ws_CUL.Range("M122").Value = -7121.5
$ws_CUL.Range("N122").Value = -25769.2
# --- GSM ---
# Row 80
$ws_GSM.Range("H80").Value = 95775.586
$ws_GSM.Range("I80").Value = 161258.42
$ws_GSM.Range("J80").Value = 4099.6
$ws_GSM.Range("K80").Value = 161258.42
$ws_GSM.Range("L80").Value = 4099.6
$ws_GSM.Range("M80").Value = -160260.42
$ws_GSM.Range("N80").Value = -6095.6
# Row 83
$ws_GSM.Range("H83").Value = 95775.586
$ws_GSM.Range("I83").Value = 161258.42
$ws_GSM.Range("J83").Value = 4099.6
$ws_GSM.Range("K83").Value = 806292.1000000001
$ws_GSM.Range("L83").Value = 20498
$ws_GSM.Range("M83").Value = -801300.1000000001
$ws_GSM.Range("N83").Value = -30482
# Row 122
$ws_GSM.Range("H122").Value = 3472.8518
$ws_GSM.Range("I122").Value = 1579.2727
$ws_GSM.Range("J122").Value = 4774.6875
$ws_GSM.Range("K122").Value = 4737.8181
$ws_GSM.Range("L122").Value = 14324.0625
$ws_GSM.Range("M122").Value = -2287.8181
$ws_GSM.Range("N122").Value = -19224.0625
# --- LTW ---
# Row 7
$ws_LTW.Range("H7").Value = 4060.2856
$ws_LTW.Range("I7").Value = 2611.1
$ws_LTW.Range("J7").Value = 5377.727
$ws_LTW.Range("K7").Value = 2611.1
$ws_LTW.Range("L7").Value = 5377.727
$ws_LTW.Range("M7").Value = -2499.1
$ws_LTW.Range("N7").Value = -5601.727
# Row 46
$ws_LTW.Range("H46").Value = 1218.3636
$ws_LTW.Range("I46").Value = 1337.1428
$ws_LTW.Range("J46").Value = 1010.5
$ws_LTW.Range("K46").Value = 1337.1428
$ws_LTW.Range("L46").Value = 1010.5
$ws_LTW.Range("M46").Value = -1149.1428
$ws_LTW.Range("N46").Value = -1386.5
# Row 82
$ws_LTW.Range("H82").Value = 5718.5293
$ws_LTW.Range("I82").Value = 3481.5
$ws_LTW.Range("K82").Value = 3481.5
$ws_LTW.Range("M82").Value = -3120.5
# Row 85
$ws_LTW.Range("H85").Value = 5718.5293
$ws_LTW.Range("I85").Value = 3481.5
$ws_LTW.Range("K85").Value = 3481.5
$ws_LTW.Range("M85").Value = -2233.5
# Row 122
$ws_LTW.Range("H122").Value = 6584.778
$ws_LTW.Range("I122").Value = 2804
$ws_LTW.Range("J122").Value = 7665
$ws_LTW.Range("K122").Value = 8412
$ws_LTW.Range("L122").Value = 22995
$ws_LTW.Range("M122").Value = -5962
$ws_LTW.Range("N122").Value = -27895
# Row 126
$ws_LTW.Range("H126").Value = 4060.2856
$ws_LTW.Range("I126").Value = 2611.1
$ws_LTW.Range("J126").Value = 5377.727
$ws_LTW.Range("K126").Value = 7833.299999999999
$ws_LTW.Range("L126").Value = 16133.181
$ws_LTW.Range("M126").Value = -5363.299999999999
$ws_LTW.Range("N126").Value = -21073.181
# --- WVR ---
# Row 122
$ws_WVR.Range("H122").Value = 2895.311
$ws_WVR.Range("I122").Value = 2510.3125
$ws_WVR.Range("K122").Value = 7530.9375
$ws_WVR.Range("M122").Value = -5080.9375
